# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right after "总计" (shifts the other
#    quarter sheets down by one position).
# 2) Populate "2022-Q3" with the fund holdings table.
# 3) Insert a new row into "总计" for the 2022-Q3 summary figures, pushing
#    the existing quarters down one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)

# --- 1) Add the new "2022-Q3" sheet right after "总计" ----------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# --- 2) Populate "2022-Q3" ---------------------------------------------------
# Header row (copy the "no style / s=2 header" look from the 总计 sheet's own
# header cells so the new sheet matches the sibling quarter sheets exactly).
$total.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Reference "no style" cell (plain data cell, e.g. a B column entry on the
# 总计 sheet) used to strip formatting after forcing text storage below.
$plain = $total.Range("B2")
# Reference "A column" styled cell (bold/centered/bordered, s=2) used to
# stamp the index column (A) of the new sheet the same way every sibling
# quarter sheet stamps it.
$aStyle = $total.Range("A2")

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $plain.Copy()
    $cell.PasteSpecial(-4122)
}

$q3Data = @(
    @{A=0; B="001411"; C="诺安创新驱动灵活配置混合A";   D="3.98"; E="80.56"; F="2.58"; G="0.1027"; H=10},
    @{A=1; B="159851"; C="华宝中证金融科技主题ETF";     D="1.94"; E="98.27"; F="4.83"; G="0.0937"; H=4},
    @{A=2; B="002051"; C="诺安创新驱动灵活配置混合C";   D="1.80"; E="80.56"; F="2.58"; G="0.0464"; H=10},
    @{A=3; B="560660"; C="新华中证云计算50ETF";         D="1.20"; E="97.03"; F="2.46"; G="0.0295"; H=10},
    @{A=4; B="516100"; C="华夏中证金融科技主题ETF";     D="0.51"; E="96.79"; F="4.79"; G="0.0244"; H=4},
    @{A=5; B="516860"; C="博时中证金融科技主题ETF";     D="0.34"; E="98.57"; F="4.88"; G="0.0166"; H=4}
)

$r = 2
foreach ($row in $q3Data) {
    $aStyle.Copy()
    $q3.Cells.Item($r, 1).PasteSpecial(-4122)
    $q3.Cells.Item($r, 1).Value = $row.A
    Set-TextCell $q3.Cells.Item($r, 2) $row.B
    Set-TextCell $q3.Cells.Item($r, 3) $row.C
    Set-TextCell $q3.Cells.Item($r, 4) $row.D
    Set-TextCell $q3.Cells.Item($r, 5) $row.E
    Set-TextCell $q3.Cells.Item($r, 6) $row.F
    Set-TextCell $q3.Cells.Item($r, 7) $row.G
    $q3.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# --- 3) Shift "总计" rows down one and insert the 2022-Q3 summary row ------
# Work bottom-up so we don't clobber data we still need to read.
$total.Range("A7").Value = 5
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)
$total.Range("A7").Value = 5
$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 11
$total.Range("D7").Value = 0.26

$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 5
$total.Range("D6").Value = 0.22

$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.1

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.27

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.14

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.31

# --- Restore original active sheet (last tab stays selected) --------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()
